## Fruta / hortaliza, semanal
## Insert the new weekly price observation at the top of the Zanahoria log
## (row 197), pushing all subsequent rows down by one. The last existing
## row (former 312) ends up duplicated at the new row 313 by the insert,
## which matches the source data (a new week was appended while the feed
## keeps the most-recent record at the top of this block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 197:312 down to 198:313, preserving formatting (style s="2" on
# column D, the date column) the same way Excel's own Insert does.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new observation.
$ws.Range("A197").Value = 10
$ws.Range("B197").Value = "Vega Modelo de Temuco"
$ws.Range("C197").Value = "La Araucanía"
$ws.Range("D197").Value = 44719
$ws.Range("D197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E197").Value = 9
$ws.Range("F197").Value = 100114013
$ws.Range("G197").Value = "Zanahoria"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 60
$ws.Range("K197").Value = 6000
$ws.Range("L197").Value = 6000
$ws.Range("M197").Value = 6000
$ws.Range("N197").Value = "$/saco 25 kilos"
$ws.Range("O197").Value = "Región de La Araucanía"
$ws.Range("P197").Value = 240
$ws.Range("Q197").Value = 25
$ws.Range("R197").Value = "Hortaliza"
